# Update template excel list question
# - Rename the header cell B1 from "CONTENT" to "DESCRIPTION"
# - Move the active cell selection to B6 (matches the new cursor position left by the author)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("data")

$ws.Range("B1").Value = "DESCRIPTION"

$ws.Range("B6").Select()
